$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old row 3 (CO 3013696384), pushing it to row 4.
$ws.Rows.Item(3).Insert()
# Fully clear the newly inserted row so it carries no cell data/formatting
# (matches source: row 3 is completely absent from the sheet data).
$ws.Rows.Item(3).Clear()

# Append the new CO numbers in column A for rows 5-12.
$ws.Range("A5").Value = "3013696547"
$ws.Range("A6").Value = "3013696548"
$ws.Range("A7").Value = "3013696549"
$ws.Range("A8").Value = "3013696550"
$ws.Range("A9").Value = "3013696551"
$ws.Range("A10").Value = "3013696552"
$ws.Range("A11").Value = "3013696553"
$ws.Range("A12").Value = "3013696554"

# Match the author's final selection state (whole row 3 selected).
[void]$ws.Rows.Item(3).Select()
